# Update LEM-to-LED mapping sheet with four new rows (96-99)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 96
$ws.Range("A96").Value = "LEM-198-32-3022KH"
$ws.Range("A96").Style = "Normal"
$ws.Range("B96").Value = "LED-198-H70-3022"
$ws.Range("B96").NumberFormat = "#,##0.0000"
$ws.Range("C96").Value = 32.889
$ws.Range("D96").Value = 30.98
$ws.Range("E96").Value = 75
$ws.Range("E96").NumberFormat = "#,##0.0000"

# Row 97
$ws.Range("A97").Value = "LEM-307-00-40KS"
$ws.Range("A97").Style = "Normal"
$ws.Range("B97").Value = "LED-307-S00-40"
$ws.Range("B97").NumberFormat = "#,##0.0000"
$ws.Range("C97").Value = 4.1327
$ws.Range("D97").Value = 2.8356
$ws.Range("E97").Value = 20
$ws.Range("E97").NumberFormat = "#,##0.0000"

# Row 98
$ws.Range("A98").Value = "LEM-313-00-2722KH"
$ws.Range("A98").Style = "Normal"
$ws.Range("B98").Value = "LED-313-H00-2722"
$ws.Range("B98").NumberFormat = "#,##0.0000"
$ws.Range("C98").Value = 25.925
$ws.Range("D98").Value = 25.3
$ws.Range("E98").Value = 65
$ws.Range("E98").NumberFormat = "#,##0.0000"

# Row 99
$ws.Range("A99").Value = "LEM-401-00-3018KH"
$ws.Range("A99").Style = "Normal"
$ws.Range("B99").Value = "LED-401-H00-3018"
$ws.Range("B99").NumberFormat = "#,##0.0000"
$ws.Range("C99").Value = 6.0401
$ws.Range("D99").Value = 4.743
$ws.Range("E99").Value = 20
$ws.Range("E99").NumberFormat = "#,##0.0000"

# Update selection / active cell to match final state
$ws.Range("A96:E99").Select() | Out-Null
